# Apply updated TPM-derived values to the LR-pairs worksheet (Fgf9-Fgfr3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.077011333333332
$ws.Range("H2").Value = 12.231034
$ws.Range("I2").Value = 0.9715624748044627
$ws.Range("J2").Value = 0.9715624748044628
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 15.30076804128977
$ws.Range("R2").Value = 137.706912371608
$ws.Range("S2").Value = 0.6660664371320496
$ws.Range("T2").Value = 0.6660664371320496

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.077011333333332
$ws.Range("H3").Value = 12.231034
$ws.Range("I3").Value = 0.9715624748044627
$ws.Range("J3").Value = 0.9715624748044628
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 5.285631830073554
$ws.Range("R3").Value = 47.57068647066199
$ws.Range("S3").Value = 0.2300918458177006
$ws.Range("T3").Value = 0.2300918458177006

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.077011333333332
$ws.Range("H4").Value = 12.231034
$ws.Range("I4").Value = 0.9715624748044627
$ws.Range("J4").Value = 0.9715624748044628
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 1.732172625117777
$ws.Range("R4").Value = 15.58955362606
$ws.Range("S4").Value = 0.07540419185471248
$ws.Range("T4").Value = 0.07540419185471249

# Row 5
$ws.Range("I5").Value = 0.02843752519553723
$ws.Range("J5").Value = 0.02843752519553723
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 0.4478517727568889
$ws.Range("R5").Value = 4.030665954812
$ws.Range("S5").Value = 0.01949569027113415
$ws.Range("T5").Value = 0.01949569027113415

# Row 6
$ws.Range("I6").Value = 0.02843752519553723
$ws.Range("J6").Value = 0.02843752519553723
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("S6").Value = 0.006734762645135534
$ws.Range("T6").Value = 0.006734762645135535

# Row 7
$ws.Range("I7").Value = 0.02843752519553723
$ws.Range("J7").Value = 0.02843752519553723
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 0.0507004993988889
$ws.Range("R7").Value = 0.4563044945900001
$ws.Range("S7").Value = 0.002207072279267553
$ws.Range("T7").Value = 0.002207072279267553
